# The commit italicizes the four verse excerpts (the quoted lines from the
# poem) that appear in the "Answer" section, and restores the leading tab
# on the "The people at the town ..." paragraph so it matches its sibling
# paragraphs (all of which start with a tab character).

$d = $word.ActiveDocument

function Italicize-Phrase {
    param([string]$Phrase)
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($Phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Italic = $true
    }
    [void]$found
}

# Paragraph: "The poem at the beginning highlights the ball ..."
Italicize-Phrase "rose with its voluptuous swell"
Italicize-Phrase "all went merry as a marriage bell"
Italicize-Phrase "nearer, clearer, deadlier than before"

# Paragraph: "The people at the town became pale with fear ..." is missing
# the leading tab that every other body paragraph has — add it back.
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("The people at the town became pale", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insertRange = $d.Range($rng2.Start, $rng2.Start)
    [void]$insertRange.InsertBefore([char]9)
}

Italicize-Phrase "trodden like the grass"

# Paragraph: "Last stanza of the poem reflects on the transitoriness ..."
Italicize-Phrase "one red  burial blent"

Write-Output "done"
